$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update start date (A2) - must remain text, not auto-converted to a date serial.
# Leading apostrophe forces Excel to treat the value as text.
$ws.Range("A2").Value = "'2020-01-02"

# Update numeric performance figures to reflect fee and tax adjustments
$ws.Range("A6").Value = 2.65
$ws.Range("A7").Value = 0.23
$ws.Range("A8").Value = 1.64
$ws.Range("A10").Value = 2.38
$ws.Range("A11").Value = 1.68
$ws.Range("A12").Value = 1.35
$ws.Range("A14").Value = 186
$ws.Range("A15").Value = 0.35
$ws.Range("A16").Value = 1.33
$ws.Range("A17").Value = 1.01
$ws.Range("A18").Value = 1.35
$ws.Range("A19").Value = 1.57
$ws.Range("A20").Value = 0.78
$ws.Range("A21").Value = 1.17
$ws.Range("A23").Value = 4.43
$ws.Range("A30").Value = 0.23
$ws.Range("A31").Value = 0.23
$ws.Range("A32").Value = 0.23
$ws.Range("A34").Value = 27
$ws.Range("A35").Value = 4.3
$ws.Range("A36").Value = 0.08
$ws.Range("A37").Value = 1.53
